# Auto-update stock values: 2025-12-11 07:55:54 UTC
#
# Appends a new trading-day column (2025-12-10) to the end of every data
# sheet in the workbook (all sheets except the first "종목" sheet), copying
# the existing header cell's formatting onto the new header cell and
# filling in the three data rows for that date.
#
# ColumnWidth values of 9.14 / 11.14 are used (instead of 10 / 12) because
# Excel's ColumnWidth property is expressed in "characters" and gets
# re-quantized internally; these inputs are the values that round-trip
# back to the workbook's native column-width units of exactly 10 and 12,
# matching the existing columns on each sheet.

$wb = $excel.ActiveWorkbook

function Add-DailyColumn {
    param($SheetIndex, $ColWidth, $HeaderValue, $Row2Value, $Row3Value)

    $ws = $wb.Worksheets.Item($SheetIndex)

    # Last used column (header row) is one less than the new column we add.
    $lastCol = $ws.UsedRange.Columns.Count
    $newCol = $lastCol + 1

    # Match the existing column width of the sheet's data columns.
    $ws.Columns.Item($newCol).ColumnWidth = $ColWidth

    # Copy the header cell's style (bold font + gray fill) onto the new
    # header cell, then set its value.
    $srcHeader = $ws.Cells.Item(1, $lastCol)
    $dstHeader = $ws.Cells.Item(1, $newCol)
    $srcHeader.Copy()
    $dstHeader.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $dstHeader.Value = $HeaderValue

    $ws.Cells.Item(2, $newCol).Value = $Row2Value
    $ws.Cells.Item(3, $newCol).Value = $Row3Value
}

$newDate = 20251210

# 시가 (open)
Add-DailyColumn 2 11.14 $newDate 623.85 55.65

# 고가 (high)
Add-DailyColumn 3 11.14 $newDate 629.21 57.09

# 저가 (low)
Add-DailyColumn 4 11.14 $newDate 620.99 54.89

# 종가 (close)
Add-DailyColumn 5 11.14 $newDate 627.61 56.65

# 거래량 (volume)
Add-DailyColumn 6 11.14 $newDate 55031384 100261104

# s20
Add-DailyColumn 7 9.14 $newDate 100 16

# s60
Add-DailyColumn 8 9.14 $newDate 84 14

# z20
Add-DailyColumn 9 9.14 $newDate 60 -25

# z60
Add-DailyColumn 10 9.14 $newDate 70 -81

# gap
Add-DailyColumn 11 11.14 $newDate 102 83

# std
Add-DailyColumn 12 11.14 $newDate 5.56 22.43

# quant
Add-DailyColumn 13 9.14 $newDate 47 73

Write-Output "Daily update column added to all data sheets."
